$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in column F (rows 2-7) forward by 14 days (2 weeks),
# keeping the existing date formatting/style intact.
$ws.Range("F2").Value = [DateTime]::FromOADate(44634)
$ws.Range("F3").Value = [DateTime]::FromOADate(44633)
$ws.Range("F4").Value = [DateTime]::FromOADate(44632)
$ws.Range("F5").Value = [DateTime]::FromOADate(44631)
$ws.Range("F6").Value = [DateTime]::FromOADate(44630)
$ws.Range("F7").Value = [DateTime]::FromOADate(44629)
